# Pós reunião Gustavo Código funcionando!
$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# ---------------------------------------------------------------------
# 1) Row 15 edits on Plan1
#    B15: "1.0" -> "50.0"   (kept as text, reuses existing shared string)
#    D15: 100000 -> 10000
#    E15: "10.0" -> "10.0 / 0.0"  (new text)
#    L15: new cell "Passando valores para ouro lado de novo"
# ---------------------------------------------------------------------

# Write the two brand-new strings first (L15 then E15) so they land in the
# shared-string table in the same order the original author typed them.
$ws1.Range("L15").Value = "Passando valores para ouro lado de novo"
$ws1.Range("E15").Value = "10.0 / 0.0"

$ws1.Range("D15").Value = 10000

# B15 needs to stay a *text* cell (matches the "50.0" text already used by
# B13/B14) instead of being auto-converted to the number 50, so force a
# text format, assign it, then pick the formatting back up from a sibling
# cell that already has the correct look (center aligned, General format).
$ws1.Range("B15").NumberFormat = "@"
$ws1.Range("B15").Value = "50.0"
$ws1.Range("B13").Copy()
$ws1.Range("B15").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# ---------------------------------------------------------------------
# 2) New cell Q8 = 24.33 (extends the used range to A1:Q15)
# ---------------------------------------------------------------------
$ws1.Range("Q8").Value = 24.33

# ---------------------------------------------------------------------
# 3) Selection moves to G14 on Plan1
# ---------------------------------------------------------------------
$ws1.Range("G14").Select()

# ---------------------------------------------------------------------
# 4) New worksheet "Duvidas" after Plan1
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $ws1)
$ws2.Name = "Duvidas"

$ws2.Range("A1").Value = "duvidas"
$ws2.Range("A2").Value = 1
$ws2.Range("B2").Value = "Equacionamento"
$ws2.Range("C2").Value = "Tenho que multiplicar as ccs por dt?"
$ws2.Range("A3").Value = 2
$ws2.Range("A4").Value = 3

$ws2.Columns.Item(2).ColumnWidth = 16
$ws2.Columns.Item(3).ColumnWidth = 33.140625

$ws2.Range("C3").Select()

$ws1.Select()
